# Auto-generated Excel COM-interop script applying the weekly crime-stat refresh
# (Volume/week-of header bump + revised precinct figures for rows 14-31).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text: volume "Number" 51 -> 52, report week 12/15-12/21 -> 12/22-12/28
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "52"
$ws.Range("C9").Characters(27, 10).Text = "12/22/2025"
$ws.Range("C9").Characters(48, 10).Text = "12/28/2025"

# ---------------------------------------------------------------------------
# 2) Cells that flip from a numeric entry to the "no data" placeholder text
#    (shared strings "0" / "***.*", using the existing label cells C14/E14 as
#    format+value templates so the right hand style (13) comes along for free).
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F22").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("G31").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H31").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Refreshed numeric figures (counts + recomputed % changes) for rows 14-31
# ---------------------------------------------------------------------------
# Row 15
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 21
$ws.Range("K15").Value = 16.666666666666
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = -4.545454545454
$ws.Range("N15").Value = -76.404494382022
# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -18.181818181818
$ws.Range("J16").Value = 226
$ws.Range("K16").Value = 3.982300884955
$ws.Range("L16").Value = 10.849056603773
$ws.Range("M16").Value = -16.961130742049
$ws.Range("N16").Value = -82.019892884468
# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -8.333333333333
$ws.Range("I17").Value = 314
$ws.Range("J17").Value = 347
$ws.Range("K17").Value = -9.510086455331
$ws.Range("L17").Value = -8.454810495626
$ws.Range("M17").Value = 41.441441441441
$ws.Range("N17").Value = -62.304921968787
# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 142
$ws.Range("J18").Value = 124
$ws.Range("K18").Value = 14.516129032258
$ws.Range("L18").Value = -5.333333333333
$ws.Range("M18").Value = -24.064171122994
$ws.Range("N18").Value = -92.577104025091
# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -6.976744186046
$ws.Range("I19").Value = 578
$ws.Range("J19").Value = 633
$ws.Range("K19").Value = -8.688783570300
$ws.Range("L19").Value = 2.664298401420
$ws.Range("M19").Value = 49.740932642487
$ws.Range("N19").Value = -51.016949152542
# Row 20
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 195
$ws.Range("J20").Value = 166
$ws.Range("K20").Value = 17.469879518072
$ws.Range("L20").Value = -13.333333333333
$ws.Range("M20").Value = 114.285714285714
$ws.Range("N20").Value = -89.291598023064
# Row 21
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = -7.692307692307
$ws.Range("F21").Value = 88
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -14.563106796116
$ws.Range("I21").Value = 1488
$ws.Range("J21").Value = 1520
$ws.Range("K21").Value = -2.105263157894
$ws.Range("L21").Value = -1.522170747849
$ws.Range("M21").Value = 24.518828451882
$ws.Range("N21").Value = -79.313221187265
# Row 22
$ws.Range("H22").Value = -100
$ws.Range("L22").Value = -54.761904761904
$ws.Range("M22").Value = -29.629629629629
# Row 23
$ws.Range("G23").Value = 1
$ws.Range("L23").Value = -42.424242424242
$ws.Range("M23").Value = -17.391304347826
# Row 24
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 94
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = -5.050505050505
$ws.Range("I24").Value = 1375
$ws.Range("J24").Value = 1447
$ws.Range("K24").Value = -4.975812024879
$ws.Range("L24").Value = -3.032440056417
$ws.Range("M24").Value = 121.061093247588
# Row 25
$ws.Range("C25").Value = 7
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = -23.913043478260
$ws.Range("I25").Value = 587
$ws.Range("J25").Value = 688
$ws.Range("K25").Value = -14.680232558139
$ws.Range("L25").Value = -8.850931677018
# Row 26
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 23.076923076923
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = 42.857142857142
$ws.Range("I26").Value = 582
$ws.Range("J26").Value = 568
$ws.Range("K26").Value = 2.464788732394
$ws.Range("L26").Value = 10.227272727272
$ws.Range("M26").Value = 6.398537477148
# Row 27
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 27
$ws.Range("K27").Value = -3.571428571428
$ws.Range("L27").Value = 8
# Row 28
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 54
$ws.Range("J28").Value = 62
$ws.Range("K28").Value = -12.903225806451
$ws.Range("L28").Value = -3.571428571428
# Row 29
$ws.Range("H29").Value = -100
$ws.Range("L29").Value = -61.904761904761
$ws.Range("N29").Value = -95.505617977528
# Row 30
$ws.Range("H30").Value = -100
$ws.Range("L30").Value = -50
$ws.Range("N30").Value = -94.771241830065
